# Update column F ("dSF") values to match re-pulled / re-pushed data and
# recalculated means for the fulmer_carson.xlsx dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -7
    "F4"  = -6
    "F5"  = -1
    "F6"  = -6
    "F7"  = -2
    "F8"  = -3
    "F9"  = -5
    "F10" = 6
    "F12" = 2
    "F15" = -6
    "F20" = -6
    "F22" = -10
    "F25" = -6
    "F26" = -1
    "F27" = 1
    "F34" = -10
    "F35" = -2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
